# Update the "取得日時" (retrieved timestamp) column for the lancers
# job-listing sheet to reflect a new fetch run at 2025-10-25 18:29:48.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-25 18:29:48"

for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
